$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sp2")

# Week 1 - Bryan - "Improve input and output UI"? No: row 33 is under
# backlog item "Add option for inputing concrete temp" (row 32), task
# "Estimated" (row 33, person Bryan). Fill in hours worked up to Dec 6
# (column I = 12/6/2014).
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = -0.5

# Week 2 - Bryan - under backlog item "Improve input and output UI"
# (row 40), task "Estimated" (row 41, person Bryan).
$ws.Range("E41").Value = 2
$ws.Range("F41").Value = 2
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 2
$ws.Range("I41").Value = -3.5
